# Commit: "add dynamic text in wat zie ik? for opleiding ouders"
#
# The "outcome" worksheet contains rows describing outcome variables. Two
# pieces of descriptive text need to be updated:
#   1. The "population" text "leerlingen groep 8" becomes
#      "leerlingen in groep 8" for the rows that reference groep-8 pupils
#      (column F).
#   2. The "type" / category text "Geld" becomes "Werk en inkomen" for the
#      rows in the work/income block (column D).
# The sheet's scroll position is also nudged down two rows, matching the
# author's view state when they made the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outcome")

# 1) Column F: "leerlingen groep 8" -> "leerlingen in groep 8"
$popRows = @(5, 6, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 33)
foreach ($r in $popRows) {
    $ws.Cells.Item($r, 6).Value = "leerlingen in groep 8"
}

# 2) Column D: "Geld" -> "Werk en inkomen"
$typeRows = @(36, 37, 38, 39, 40, 41, 42, 43, 44)
foreach ($r in $typeRows) {
    $ws.Cells.Item($r, 4).Value = "Werk en inkomen"
}

# 3) Update the sheet's view so the visible top-left cell matches the
#    author's saved scroll position.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
